$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row additions: Year 2024 (V1) and Year 2025 (W1) ---
$ws.Range("V1").Value = "Year 2024"
$ws.Range("W1").Value = "Year 2025"

# --- Column V ("Year 2024") gets a fresh run of month-start dates (Jul-2023..Jun-2024) ---
# Match the existing date format used throughout column A-U (style index 1 / numFmtId 17 "mmm-yy")
# Values are plain Excel date-serial numbers (days since 1899-12-30), same as the rest of the sheet.
$ws.Range("V2:V13").NumberFormat = "mmm-yy"
$ws.Range("V2").Value = 45108
$ws.Range("V3").Value = 45139
$ws.Range("V4").Value = 45170
$ws.Range("V5").Value = 45200
$ws.Range("V6").Value = 45231
$ws.Range("V7").Value = 45261
$ws.Range("V8").Value = 45292
$ws.Range("V9").Value = 45323
$ws.Range("V10").Value = 45352
$ws.Range("V11").Value = 45383
$ws.Range("V12").Value = 45413
$ws.Range("V13").Value = 45444

# --- The two dates that moved off column U (Year 2023) into V get cleared from U,
#     but keep the "not-yet-reached" placeholder look of column T (format-only copy) ---
$ws.Range("T3").Copy()
$ws.Range("U3").PasteSpecial(-4122)
$ws.Range("U3").Value = ""

$ws.Range("T4").Copy()
$ws.Range("U4").PasteSpecial(-4122)
$ws.Range("U4").Value = ""

# --- New "Fixed"/"Variable" mini cost table ---
$ws.Range("D21").Value = "Fixed"
$ws.Range("E21").Value = "Variable"

$ws.Range("D22").Value = 3309250
$ws.Range("E22").Formula = "=20%*D22"
$ws.Range("G22").Formula = "=D22+E22"
$ws.Range("H22").Formula = "=G22*20%"

$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 222513
$ws.Range("E23").Formula = "=D22+E22"
$ws.Range("F23").NumberFormat = "#,##0"
$ws.Range("F23").Formula = "=E23+D23"
$ws.Range("H23").Formula = "=G22+H22"

$ws.Range("D24").Formula = "=SUM(D22:D23)"
$ws.Range("H24").Formula = "=D24*20%"

$ws.Range("H25").Formula = "=D24+H24"
$ws.Range("I25").Formula = "=H25*20%"

$ws.Range("I26").Formula = "=H25+I25"

# --- Restore selection cursor to match the saved workbook state ---
[void]$ws.Range("J20").Select()
